$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D2").Value = "25.866.56"
$ws.Range("E2").Value = "  -0.48%  "

$ws.Range("D3").Value = "1.619.71"
$ws.Range("E3").Value = "  -1.22%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'212.74"
$ws.Range("E5").Value = "  -1.15%  "

$ws.Range("D6").Value = "'0.499"
$ws.Range("E6").Value = "  -1.35%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  -1.35%  "

$ws.Range("D9").Value = "'0.0617"
$ws.Range("E9").Value = "  -2.97%  "

$ws.Range("D10").Value = "'18.46"
$ws.Range("E10").Value = "  -4.94%  "

$ws.Range("E11").Value = "  -0.63%  "

$ws.Range("D12").Value = "1.844.48"
$ws.Range("E12").Value = "  -1.20%  "

$ws.Range("D13").Value = "1.616.61"
$ws.Range("E13").Value = "  -1.80%  "


$ws.Range("E15").Value = "  -3.28%  "

$ws.Range("D16").Value = "25.871.03"
$ws.Range("E16").Value = "  -0.47%  "

$ws.Range("D17").Value = "'61.50"
$ws.Range("E17").Value = "  -2.55%  "

$ws.Range("E18").Value = "  -2.81%  "

$ws.Range("D19").Value = "'1.00"
$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("D20").Value = "'191.57"
$ws.Range("E20").Value = "  -0.97%  "

$ws.Range("E21").Value = "  -1.65%  "

$ws.Range("D22").Value = "'9.47"
$ws.Range("E22").Value = "  -2.62%  "

$ws.Range("E23").Value = "  -2.12%  "

$ws.Range("E24").Value = "  +2.20%  "

$ws.Range("D25").Value = "'143.78"
$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  -3.62%  "

$ws.Range("E28").Value = "  -2.40%  "

$ws.Range("D29").Value = "'15.15"
$ws.Range("E29").Value = "  -1.97%  "

$ws.Range("E30").Value = "  -1.21%  "

$ws.Range("E31").Value = "  -2.46%  "

$ws.Range("D32").Value = "'3.13"
$ws.Range("E32").Value = "  -3.78%  "

$ws.Range("E33").Value = "  -4.82%  "

$ws.Range("E34").Value = "  -1.98%  "

$ws.Range("E35").Value = "  -2.75%  "

$ws.Range("D36").Value = "1.123.65"
$ws.Range("E36").Value = "  -0.10%  "

$ws.Range("D37").Value = "'0.838"
$ws.Range("E37").Value = "  -6.58%  "

$ws.Range("E38").Value = "  -3.78%  "

$ws.Range("E39").Value = "  -3.99%  "

$ws.Range("E40").Value = "  -2.01%  "

$ws.Range("D41").Value = "'98.04"
$ws.Range("E41").Value = "  -0.23%  "

$ws.Range("D42").Value = "1.755.32"
$ws.Range("E42").Value = "  -1.02%  "

$ws.Range("E43").Value = "  -5.74%  "

$ws.Range("D44").Value = "'5.03"
$ws.Range("E44").Value = "  -5.28%  "

$ws.Range("D45").Value = "0.0₆0113"
$ws.Range("E45").Value = "  -1.26%  "

$ws.Range("E46").Value = "  +0.86%  "

$ws.Range("D47").Value = "'54.01"
$ws.Range("E47").Value = "  -3.71%  "

$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("E49").Value = "  -0.61%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.48"
$ws.Range("E50").Value = "  -2.90%  "

$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  +0.05%  "
